# Update the lattice-multiplication exercise table: every cell's
# problem (top line), partial-product digits, divider and the two
# lattice-row leading digits are replaced with new values, keeping the
# existing 5x3 table shape and per-run formatting (sz=32) intact.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$sep = [char]11   # corresponds to a <w:br/> line break inside a cell

$cellText = "97 x 50" + $sep + "  5    0" + $sep + "  ----" + $sep + "9|    |" + $sep + "7|    |"
$t.Cell(1, 1).Range.Text = $cellText
$cellText = "61 x 83" + $sep + "  8    3" + $sep + "  ----" + $sep + "6|    |" + $sep + "1|    |"
$t.Cell(1, 2).Range.Text = $cellText
$cellText = "72 x 35" + $sep + "  3    5" + $sep + "  ----" + $sep + "7|    |" + $sep + "2|    |"
$t.Cell(1, 3).Range.Text = $cellText

$cellText = "84 x 88" + $sep + "  8    8" + $sep + "  ----" + $sep + "8|    |" + $sep + "4|    |"
$t.Cell(2, 1).Range.Text = $cellText
$cellText = "98 x 91" + $sep + "  9    1" + $sep + "  ----" + $sep + "9|    |" + $sep + "8|    |"
$t.Cell(2, 2).Range.Text = $cellText
$cellText = "19 x 91" + $sep + "  9    1" + $sep + "  ----" + $sep + "1|    |" + $sep + "9|    |"
$t.Cell(2, 3).Range.Text = $cellText

$cellText = "43 x 39" + $sep + "  3    9" + $sep + "  ----" + $sep + "4|    |" + $sep + "3|    |"
$t.Cell(3, 1).Range.Text = $cellText
$cellText = "60 x 27" + $sep + "  2    7" + $sep + "  ----" + $sep + "6|    |" + $sep + "0|    |"
$t.Cell(3, 2).Range.Text = $cellText
$cellText = "28 x 26" + $sep + "  2    6" + $sep + "  ----" + $sep + "2|    |" + $sep + "8|    |"
$t.Cell(3, 3).Range.Text = $cellText

$cellText = "92 x 97" + $sep + "  9    7" + $sep + "  ----" + $sep + "9|    |" + $sep + "2|    |"
$t.Cell(4, 1).Range.Text = $cellText
$cellText = "13 x 67" + $sep + "  6    7" + $sep + "  ----" + $sep + "1|    |" + $sep + "3|    |"
$t.Cell(4, 2).Range.Text = $cellText
$cellText = "93 x 29" + $sep + "  2    9" + $sep + "  ----" + $sep + "9|    |" + $sep + "3|    |"
$t.Cell(4, 3).Range.Text = $cellText

$cellText = "33 x 57" + $sep + "  5    7" + $sep + "  ----" + $sep + "3|    |" + $sep + "3|    |"
$t.Cell(5, 1).Range.Text = $cellText
$cellText = "20 x 13" + $sep + "  1    3" + $sep + "  ----" + $sep + "2|    |" + $sep + "0|    |"
$t.Cell(5, 2).Range.Text = $cellText
$cellText = "68 x 48" + $sep + "  4    8" + $sep + "  ----" + $sep + "6|    |" + $sep + "8|    |"
$t.Cell(5, 3).Range.Text = $cellText
